$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> FAPs)
$ws.Range("I2").Value = 0.2503572190582515
$ws.Range("J2").Value = 0.2503572190582515
$ws.Range("M2").Value = 0.6946430000000001
$ws.Range("N2").Value = 2.083929
$ws.Range("O2").Value = 0.1140293552421611
$ws.Range("P2").Value = 0.1140293552421611
$ws.Range("Q2").Value = 1.403846643470667
$ws.Range("S2").Value = 0.02854807226943292
$ws.Range("T2").Value = 0.02854807226943292

# Row 3 (FAPs -> MuSCs)
$ws.Range("I3").Value = 0.2503572190582515
$ws.Range("J3").Value = 0.2503572190582515
$ws.Range("M3").Value = 5.373609333333333
$ws.Range("N3").Value = 16.120828
$ws.Range("O3").Value = 0.8821066470161785
$ws.Range("P3").Value = 0.8821066470161785
$ws.Range("Q3").Value = 10.85985668310578
$ws.Range("R3").Value = 97.73871014795201
$ws.Range("S3").Value = 0.2208417670597692
$ws.Range("T3").Value = 0.2208417670597692

# Row 4 (FAPs -> Resolving-Mac)
$ws.Range("I4").Value = 0.2503572190582515
$ws.Range("J4").Value = 0.2503572190582515
$ws.Range("M4").Value = 0.02353866666666667
$ws.Range("N4").Value = 0.070616
$ws.Range("O4").Value = 0.00386399774166032
$ws.Range("P4").Value = 0.00386399774166032
$ws.Range("Q4").Value = 0.04757073517155556
$ws.Range("R4").Value = 0.428136616544
$ws.Range("S4").Value = 0.000967379729049442
$ws.Range("T4").Value = 0.000967379729049442

# Row 5 (MuSCs -> FAPs)
$ws.Range("G5").Value = 6.051349666666667
$ws.Range("H5").Value = 18.154049
$ws.Range("I5").Value = 0.7496427809417484
$ws.Range("J5").Value = 0.7496427809417485
$ws.Range("M5").Value = 0.6946430000000001
$ws.Range("N5").Value = 2.083929
$ws.Range("O5").Value = 0.1140293552421611
$ws.Range("P5").Value = 0.1140293552421611
$ws.Range("Q5").Value = 4.203527686502334
$ws.Range("R5").Value = 37.83174917852101
$ws.Range("S5").Value = 0.08548128297272821
$ws.Range("T5").Value = 0.08548128297272822

# Row 6 (MuSCs -> MuSCs)
$ws.Range("G6").Value = 6.051349666666667
$ws.Range("H6").Value = 18.154049
$ws.Range("I6").Value = 0.7496427809417484
$ws.Range("J6").Value = 0.7496427809417485
$ws.Range("M6").Value = 5.373609333333333
$ws.Range("N6").Value = 16.120828
$ws.Range("O6").Value = 0.8821066470161785
$ws.Range("P6").Value = 0.8821066470161785
$ws.Range("Q6").Value = 32.51758904806356
$ws.Range("R6").Value = 292.658301432572
$ws.Range("S6").Value = 0.6612648799564093
$ws.Range("T6").Value = 0.6612648799564094

# Row 7 (MuSCs -> Resolving-Mac)
$ws.Range("G7").Value = 6.051349666666667
$ws.Range("H7").Value = 18.154049
$ws.Range("I7").Value = 0.7496427809417484
$ws.Range("J7").Value = 0.7496427809417485
$ws.Range("M7").Value = 0.02353866666666667
$ws.Range("N7").Value = 0.070616
$ws.Range("O7").Value = 0.00386399774166032
$ws.Range("P7").Value = 0.00386399774166032
$ws.Range("Q7").Value = 0.1424407026871111
$ws.Range("R7").Value = 1.281966324184
$ws.Range("S7").Value = 0.002896618012610878
$ws.Range("T7").Value = 0.002896618012610878
